$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: existing company record (id changes text "1" -> "2"); metrics recomputed ---
$ws.Range("B2").Value = "'2"
$ws.Range("G2").Value = -0.8100426338228327
$ws.Range("H2").Value = -5.613453339649454
$ws.Range("I2").Value = -5.141133279328181
$ws.Range("J2").Value = -5.141133279328181
$ws.Range("K2").Value = -110.4
$ws.Range("L2").Value = -5.229748934154428
$ws.Range("U2").Value = 46.59999999999999
$ws.Range("V2").Value = 0.1081206496519721
$ws.Range("W2").Value = -5.219830295411691
$ws.Range("X2").Value = 0.06477037355041163
$ws.Range("Y2").Value = -5.284600668962102
$ws.Range("Z2").Value = 0.4335893727325336
$ws.Range("AA2").Value = 0.4517368329048101
$ws.Range("AB2").Value = 0.06226114572961206
$ws.Range("AC2").Value = 0.3894756871751981
$ws.Range("AD2").Value = 21.78
$ws.Range("AE2").Value = 0.6866176330895362
$ws.Range("AF2").Value = 22.46661763308954
$ws.Range("AG2").Value = -24.13338236691046
$ws.Range("AH2").Value = 0.04954414891741338
$ws.Range("AI2").Value = 0.3046873908006267
$ws.Range("AJ2").Value = -0.05931521860236229
$ws.Range("AK2").Value = -0.8893290495232156
$ws.Range("AL2").Value = 2.547
$ws.Range("AM2").Value = 1.28
$ws.Range("AN2").Value = -0.2112676056338028
$ws.Range("AO2").Value = -42.83470749901846
$ws.Range("AP2").Value = 0.2340955880855008
$ws.Range("AQ2").Value = -85.23437499999999

# --- Row 3: company renamed to Stealth BioTherapeutics Corp; metrics recomputed ---
$ws.Range("B3").Value = "Stealth BioTherapeutics Corp (NasdaqGM:MITO)"
$ws.Range("G3").Value = -0.5165876777251184
$ws.Range("H3").Value = -2.156398104265403
$ws.Range("I3").Value = -1.570110119744925
$ws.Range("J3").Value = -1.570110119744925
$ws.Range("K3").Value = -35.4
$ws.Range("L3").Value = -1.677725118483412
$ws.Range("U3").Value = 19.9
$ws.Range("V3").Value = 0.2525380710659899
$ws.Range("W3").Value = -9.567567567567567
$ws.Range("X3").Value = 0.0669883970030922
$ws.Range("Y3").Value = -9.634555964570659
$ws.Range("Z3").Value = -1.317647922003059
$ws.Range("AA3").Value = 2.068852336597874
$ws.Range("AB3").Value = 0.06307146927988766
$ws.Range("AC3").Value = 2.005780867317986
$ws.Range("AD3").Value = 8.98
$ws.Range("AE3").Value = 0.6866176330895362
$ws.Range("AF3").Value = 9.666617633089537
$ws.Range("AG3").Value = -10.23338236691046
$ws.Range("AH3").Value = 0.1092685341851918
$ws.Range("AI3").Value = 0.8098288753333838
$ws.Range("AJ3").Value = -0.1492472973024695
$ws.Range("AK3").Value = 1.285054753798126
$ws.Range("AL3").Value = 2.08
$ws.Range("AM3").Value = 1.726
$ws.Range("AN3").Value = -0.2746849382111832
$ws.Range("AO3").Value = -16.20192307692308
$ws.Range("AP3").Value = 0.3130240538024734
$ws.Range("AQ3").Value = -19.52491309385864

# --- Row 4: new row, Tanvex BioPharma, Inc. (TSEC:6541) added back as its own record ---
$ws.Range("A4").Value = "Cayman Islands"
$ws.Range("B4").Value = "Tanvex BioPharma, Inc. (TSEC:6541)"
$ws.Range("C4").Value = "Drugs (Biotechnology)"
$ws.Range("G4").Value = -620.0000000000002
$ws.Range("H4").Value = -7300
$ws.Range("I4").Value = -7540
$ws.Range("J4").Value = -7540
$ws.Range("K4").Value = -75
$ws.Range("L4").Value = -7500
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 26.7
$ws.Range("V4").Value = 0.075809199318569
$ws.Range("W4").Value = -0.8720930232558139
$ws.Range("X4").Value = 0.06255235009773107
$ws.Range("Y4").Value = -0.934645373353545
$ws.Range("Z4").Value = 0.0001545595054095827
$ws.Range("AA4").Value = -1.165378670788253
$ws.Range("AB4").Value = 0.06145082217933646
$ws.Range("AC4").Value = -1.22682949296759
$ws.Range("AD4").Value = 12.8
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 12.8
$ws.Range("AG4").Value = -13.9
$ws.Range("AH4").Value = 0.03506849315068494
$ws.Range("AI4").Value = 0.2071197411003236
$ws.Range("AJ4").Value = -0.04108779190067986
$ws.Range("AK4").Value = -0.396011396011396
$ws.Range("AL4").Value = 0.467
$ws.Range("AM4").Value = -0.446
$ws.Range("AN4").Value = -0.1818181818181818
$ws.Range("AO4").Value = -161.4561027837259
$ws.Range("AP4").Value = 0.1974431818181818
$ws.Range("AQ4").Value = 169.0582959641256
